$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force a run boundary at an absolute character position by
# toggling a formatting property on/off over the range running from the
# split point to the end of the whole (already-merged) segment. Toggling
# back to the original value leaves no visible trace in <w:rPr> (no
# stray <w:b/>) while still splitting the run cleanly in two. Applying
# this repeatedly, front-to-back, for a list of split points carves the
# segment into as many runs as needed.
# ---------------------------------------------------------------------
function SplitRunsAt($positions, $segEnd) {
    foreach ($pos in $positions) {
        $r = $d.Range($pos, $segEnd)
        $r.Bold = 1
        $r.Bold = 0
    }
}

# Replace a whole run/segment of text (old -> new) and then re-split the
# resulting (merged) run back into several runs according to $parts,
# which is an ordered list of the literal substrings that should each
# become their own <w:r>. The concatenation of $parts must equal $new.
function ReplaceAndSplit($old, $parts) {
    $new = [string]::Join("", $parts)
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
    $start = $rng.Start
    $segEnd = $rng.End
    $positions = New-Object System.Collections.ArrayList
    $cursor = $start
    for ($i = 0; $i -lt $parts.Length - 1; $i++) {
        $cursor = $cursor + $parts[$i].Length
        [void]$positions.Add($cursor)
    }
    SplitRunsAt $positions $segEnd
}

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Cosmic Queries: Unveiling Life's Mysteries", $false, $false, $false, $false, $false, $true, 1, $false, "Biology: Unveiling the Symphony of Life", 2) | Out-Null

# ---------------------------------------------------------------------
# Author name: "Dr. Maggie Havens" (3 runs) -> "Luna Ravenwood" (1 run)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Dr. Maggie Havens", $false, $false, $false, $false, $false, $true, 1, $false, "Luna Ravenwood", 2) | Out-Null

# ---------------------------------------------------------------------
# Email address: "maggie" -> "Luna", "havens@luminaryastronomy" -> "Ravenwood@newdawnacademy"
# Keep the surrounding "." / "org" runs intact (re-split after merge).
# ---------------------------------------------------------------------
$emailOld = "maggie.havens@luminaryastronomy.org"
$emailParts = @("Luna", ".", "Ravenwood@newdawnacademy", ".", "org")
ReplaceAndSplit $emailOld $emailParts

# ---------------------------------------------------------------------
# Main body paragraph - three "lines" separated by <w:br/>.
# Each line is its own contiguous run-segment; replacing the whole
# segment text in one shot merges it into a single run, which we then
# re-split back into the alternating sentence/period runs the target
# structure wants.
# ---------------------------------------------------------------------
$line1Old = "In the grand tapestry of the cosmos, humans, like celestial detectives, embark on an impassioned quest to unravel the cryptic script scrawled across the cosmic canvas. Our telescopes, microscopes, and analytical minds are our tools, etching inroads into the enigma of existence. We parse cosmic hieroglyphs, transcending terrestrial boundaries to explore an arena where stars venture into their final chapters, vast nebulae unfurl their spectral mantles, and distant planets whisper tales of celestial odysseys. With each discovery, we unearth clues to our origins, our evolution, and the ultimate jigsaw of the universe"
$line1Parts = @(
    "Biology - a blend of wonder, exploration, and discovery - stands as a cornerstone of scientific comprehension",
    ".",
    " We delve into the intricate realm of living organisms, unlocking mysteries that govern the very fabric of life",
    ".",
    " Within the microscopic world lies an enigmatic world teeming with biological marvels that beckon our insatiable curiosity",
    "."
)
ReplaceAndSplit $line1Old $line1Parts

$line2Old = "As we traverse the labyrinthine pathways of astronomy, we confront questions as profound as they are enigmatic: Do our cosmic origins bear the whispers of ancient alchemy, forged in celestial cauldrons? Are there other life-bearing planets circling distant stars, holding the promise of extraterrestrial encounters? Why does the universe reveal itself in its kaleidoscope of colors and patterns, urging us to decipher its cryptic composition? These enigmas, like threads of a cosmic tapestry, bind our species together in a pursuit of knowledge that transcends time, cultures, and boundaries"
$line2Parts = @(
    "Biology unravels the intricate web of life's processes, from the smallest microorganisms to the complex workings of the human body",
    ".",
    " We uncover the secrets of genetic inheritance, tracing the lineage of traits that define each individual",
    ".",
    " Moreover, we embark on quests to understand the intricate mechanisms that fuel cellular respiration, revealing how cells obtain energy from food",
    ".",
    " Our journey into biology's tapestry illuminates the processes of evolution and adaptation, highlighting the astounding diversity of living species",
    "."
)
ReplaceAndSplit $line2Old $line2Parts

$line3Old = "Our pursuit of cosmic understanding is inextricably intertwined with our introspections upon ourselves, terrestrial life, and the intricate web that connects animate and inanimate matter. As we peer into the abyssal depths of space, we necessarily embark upon an equally formidable journey into our own consciousness, unraveling the secrets hidden within the enigmatic chambers of our minds. In this celestial expedition, we are not merely spectators; we are active participants, forging a new chapter in the grand narrative of the universe, etching our indelible mark upon the cosmos itself"
$line3Parts = @(
    "Biology's tapestry intertwines the study of ecosystems and their intricate relationships",
    ".",
    " We learn of ecological interactions and the crucial role that biodiversity plays in maintaining the delicate balance of life on Earth",
    ".",
    " From microscopic organisms to towering trees, this field of science underscores interdependence in the intricate web of nature",
    ".",
    " Biology provides a profound understanding of the ecological challenges confronting our planet, inspiring us to seek solutions for sustainable living",
    "."
)
ReplaceAndSplit $line3Old $line3Parts

# ---------------------------------------------------------------------
# Summary heading body paragraph
# ---------------------------------------------------------------------
$summaryOld = "The cosmos beckons us to embark on an epoch-making expedition of exploration, inviting us to unravel the enigmatic threads of our cosmic tapestry. With fervent zeal, we seek answers to questions as profound as our own existence, simultaneously gazing outward into the enigmatic expanse of the universe and inward into the depths of our own lives. Every celestial revelation unveils a new chapter in the tome of our shared narrative, binding us eternally to the cosmos that cradles our planetary home"
$summaryParts = @(
    "Biology unveils the fascinating world of living organisms, from cellular mechanisms to ecological interactions",
    ".",
    " It illuminates the intricacies of genetic inheritance, cellular respiration, evolution, and ecosystem dynamics",
    ".",
    " Biology empowers us to appreciate the beauty and complexity of life and challenges us to address ecological issues, fostering responsible stewardship of our shared planet",
    "."
)
ReplaceAndSplit $summaryOld $summaryParts

# ---------------------------------------------------------------------
# Append a new trailing empty paragraph before the final section break.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
